$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-12-21"

# Update the header label in I1 (shared string "2022 (through 12-20)" -> "2022 (through 12-21)")
$ws.Range("I1").Value = "2022 (through 12-21)"

# Update December total-column value (I13): 84 -> 91
$ws.Range("I13").Value = 91

# Update the Total row's total-column value (I14): 1601 -> 1608
$ws.Range("I14").Value = 1608
